$wb = $excel.ActiveWorkbook

# --- 1. Update the "type of water body" translation strings -----------------
# The literal "<" character in these three prompts is replaced with the
# literal text "&#60;" (i.e. the string now contains an ampersand, hash,
# "60" and a semicolon rather than a real less-than sign).
$wsTrans = $wb.Worksheets.Item("table_specific_translations")
$wsTrans.Range("B3").Value = "53a(ii). What is the type of water body close (&#60;50m) to the household?"
$wsTrans.Range("C3").Value = "53a(ii). Qual é o tipo de corpo de água perto (&#60;50 m) do agregado familiar?"
$wsTrans.Range("D3").Value = "53a (ii). Eneo la maji lililo karibu (&#60;50m) na kaya yako ni la aina gani?"

# --- 2. Bump the form version / id on the settings sheet ---------------------
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("B3").Value = 20210305001

# --- 3. Update the remembered cursor position on table_specific_translations -
# Move to the sheet, change the selection, then return focus to the sheet
# that was originally active (settings) so the workbook-level active tab is
# left untouched.
$wsTrans.Activate()
$wsTrans.Range("D4").Select()
$wsSettings.Activate()
$wsSettings.Range("B4").Select()
